# LOBSTAHS_acyl_prop_ranges.xlsx -- "Update files with title of paper as submitted"
#
# The Notes sheet gets a new line (row 6) citing the published paper; every
# row from the old row 7 onward shifts down by one to make room for it, and
# the new citation's active cell becomes the sheet's selection.

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")

# Insert a brand-new row above the current row 6 (the blank separator row),
# which pushes rows 6.. down to 7.. -- this reproduces the row-shift seen
# throughout the diff (old r7 -> r8, r9 -> r10, r11 -> r12, etc.).
$notes.Rows.Item(6).Insert()

# Populate the freshly inserted row 6 with the paper citation.
$notes.Range("A6").Value = "See Collins, J.R., B.R. Edwards, H.F. Fredricks, and B.A.S. Van Mooy, 2016, ""LOBSTAHS: A Novel Lipidomics Strategy for Semi-Untargeted Discovery and Identification of Oxidative Stress Biomarkers"""

# Match the updated selection recorded in the saved workbook.
$notes.Range("A6").Select()
